$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 20 - new time registration entry
$ws.Range("A20").Value = "Ret OC0804 efter review "
$ws.Range("B20").Value = "System Analyst "
$ws.Range("C20").Value = 43888
$ws.Range("D20").Value = 0.4375
$ws.Range("E20").Value = 0.4513888888888889
$ws.Range("F20").Value = 0.0069444444444444441
$ws.Range("F20").NumberFormat = "h:mm"

# Row 21 - new time registration entry
$ws.Range("A21").Value = "Fælles gennemgang af CDC01 og SD01"
$ws.Range("B21").Value = "Software Architect"
$ws.Range("C21").Value = 43888
$ws.Range("D21").Value = 0.50694444444444442
$ws.Range("E21").Value = 0.66666666666666663
$ws.Range("F21").Value = "???"

# Update the active selection to match the saved workbook state
$ws.Range("D25").Select()
